$d = $word.ActiveDocument

# -----------------------------------------------------------------
# Edit 1: "Pós-Condições: ..." paragraph (Paragraphs.Item(7))
#   Before: "O usuário Gerente altera a sua atual por " + "outra" + " inserida duas vezes para confirmação,"
#   After : "O usuário" + " Gerente altera a sua senha."
# -----------------------------------------------------------------
$p1 = $d.Paragraphs.Item(7)
$r1 = $p1.Range
$start1 = $r1.Start

# Temporarily bold the middle run ("outra") so it will not be silently
# re-merged into its (identically formatted) neighbours while we edit
# the surrounding text.
$mid = $d.Range($start1 + 56, $start1 + 61)
$mid.Bold = 1

# Shrink the first run down to its new (shorter) text.
$front = $d.Range($start1 + 15, $start1 + 56)
$front.Text = "O usuário"

# The (still-bold) middle run now starts right after "O usuário".
$midLen = 5   # length of "outra"
$midNow = $d.Range($start1 + 15 + 9, $start1 + 15 + 9 + $midLen)
$midNow.Text = " Gerente altera a sua senha."

# Remove the trailing run (" inserida duas vezes para confirmação,").
$parEnd1 = $d.Paragraphs.Item(7).Range.End
$tailStart = $start1 + 15 + 9 + (" Gerente altera a sua senha.").Length
$tail = $d.Range($tailStart, $parEnd1 - 1)
$tail.Text = ""

# Restore normal (non-bold) formatting on the new run.
$parEnd1b = $d.Paragraphs.Item(7).Range.End
$final = $d.Range($start1 + 15 + 9, $parEnd1b - 1)
$final.Bold = 0

# -----------------------------------------------------------------
# Edit 2: "Clicar no campo Alterar Senha Gerência" paragraph (Paragraphs.Item(12))
#   Before: "Clicar no campo Alterar Senha Gerência" (single run)
#   After : "Clicar na opção " + "Alterar Senha Gerência" (two runs)
# -----------------------------------------------------------------
$p2 = $d.Paragraphs.Item(12)
$r2 = $p2.Range
$start2 = $r2.Start

# Temporarily bold the tail text so it survives as its own run.
$tail2 = $d.Range($start2 + 16, $start2 + 38)
$tail2.Bold = 1

# Rewrite the first 16 characters in place.
$front2 = $d.Range($start2, $start2 + 16)
$front2.Text = "Clicar na opção "

# Restore normal (non-bold) formatting on the tail run.
$parEnd2 = $d.Paragraphs.Item(12).Range.End
$tail2b = $d.Range($start2 + 16, $parEnd2 - 1)
$tail2b.Bold = 0
